# ============================================================================
# Edit: add "2022-Q1" worksheet (fund-holdings detail) before "总计" sheet,
# and add a new leading row to "总计" summarizing the 2022-Q1 quarter.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. Locate existing sheets we need as anchors / style sources.
# ----------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# ----------------------------------------------------------------------
# 2. Insert the new "2022-Q1" worksheet right after "2021-Q4" (i.e. right
#    before "总计", matching the workbook's sheet ordering in the diff).
# ----------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Re-resolve the "总计" sheet *after* inserting/positioning the new sheet:
# worksheet references obtained before a sheet-collection change can end up
# pointing at the wrong position, so look it up fresh by name now.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Match the header styling (bold font + border) used throughout the workbook
# by copying the format from an existing header cell.
$q4Sheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$fundData = @(
    @(0, '512880', '国泰中证全指证券公司ETF', '322.34', '99.85', '6.56', '21.1455', 3),
    @(1, '512000', '华宝中证全指证券公司ETF', '230.47', '99.86', '6.55', '15.0958', 3),
    @(2, '512900', '南方中证全指证券公司ETF', '78.29', '99.89', '6.56', '5.1358', 3),
    @(3, '159841', '天弘中证全指证券公司ETF', '42.72', '99.81', '6.54', '2.7939', 3),
    @(4, '512070', '易方达沪深300非银行金融ETF', '39.28', '99.37', '5.11', '2.0072', 4),
    @(5, '161720', '招商中证全指证券公司指数（LOF）', '28.99', '94.46', '6.18', '1.7916', 3),
    @(6, '001552', '天弘中证证券保险指数型发起式 A', '29.93', '95.08', '4.71', '1.4097', 4),
    @(7, '501016', '国泰中证申万证券行业指数（LOF）', '19.81', '93.96', '6.29', '1.2460', 3),
    @(8, '163113', '申万菱信中证申万证券行业指数（LOF）', '18.63', '93.17', '6.19', '1.1532', 3),
    @(9, '160633', '鹏华中证全指证券公司指数（LOF）', '17.69', '94.30', '6.17', '1.0915', 3),
    @(10, '510230', '国泰上证180金融ETF', '36.29', '99.95', '2.76', '1.0016', 8),
    @(11, '161027', '富国中证全指证券公司指数', '14.89', '94.31', '6.22', '0.9262', 3),
    @(12, '159993', '鹏华国证证券龙头ETF', '13.43', '97.76', '6.73', '0.9038', 3),
    @(13, '502010', '易方达证券公司指数（LOF）', '14.03', '94.56', '6.19', '0.8685', 3),
    @(14, '001553', '天弘中证证券保险指数型发起式 C', '18.31', '95.08', '4.71', '0.8624', 4),
    @(15, '501047', '汇添富中证全指证券公司指数LOF A', '18.61', '94.26', '4.30', '0.8002', 5),
    @(16, '515010', '华夏中证全指证券公司ETF', '11.46', '99.36', '6.50', '0.7449', 3),
    @(17, '160516', '博时中证全指证券公司指数', '8.30', '94.96', '6.19', '0.5138', 3),
    @(18, '160625', '鹏华中证800证券保险指数（LOF）', '10.45', '94.58', '4.72', '0.4932', 4),
    @(19, '501048', '汇添富中证全指证券公司指数LOF C', '10.39', '94.26', '4.30', '0.4468', 5),
    @(20, '160419', '华安中证全指证券公司指数（LOF）A', '5.41', '94.03', '6.21', '0.3360', 3),
    @(21, '159842', '银华中证全指证券公司ETF', '4.73', '97.53', '6.48', '0.3065', 3),
    @(22, '515560', '建信中证全指证券公司ETF', '4.71', '98.22', '6.01', '0.2831', 3),
    @(23, '502053', '长盛中证全指证券公司指数（LOF）', '4.28', '94.45', '6.16', '0.2636', 3),
    @(24, '013659', '中融金融鑫选3个月持有混合A', '2.83', '72.19', '8.23', '0.2329', 1),
    @(25, '515850', '富国中证全指证券公司ETF', '2.17', '99.49', '6.52', '0.1415', 3),
    @(26, '159848', '国联安中证全指证券公司ETF', '2.13', '98.30', '6.48', '0.1380', 3),
    @(27, '001304', '建信鑫安回报灵活配置混合', '2.13', '66.83', '5.60', '0.1193', 7),
    @(28, '512570', '易方达中证全指证券公司ETF', '1.74', '99.11', '6.50', '0.1131', 3),
    @(29, '013660', '中融金融鑫选3个月持有混合C', '1.15', '72.19', '8.23', '0.0946', 1),
    @(30, '399001', '中海上证50指数增强', '2.36', '93.43', '2.92', '0.0689', 9),
    @(31, '510200', '汇安上证证券ETF', '0.74', '97.70', '9.25', '0.0684', 2),
    @(32, '012605', '西藏东财中证证券保险领先指数型发起式证券投资基金A', '1.16', '94.85', '5.81', '0.0674', 5),
    @(33, '515630', '鹏华中证800证券保险ETF', '1.30', '96.69', '4.87', '0.0633', 4),
    @(34, '516980', '华富中证证券公司先锋策略ETF', '0.39', '98.42', '14.37', '0.0560', 2),
    @(35, '004836', '中融鑫价值灵活配置混合A', '0.92', '63.92', '4.77', '0.0439', 6),
    @(36, '012606', '西藏东财中证证券保险领先指数型发起式证券投资基金C', '0.58', '94.85', '5.81', '0.0337', 5),
    @(37, '510650', '华夏金融ETF', '0.73', '99.00', '3.19', '0.0233', 8),
    @(38, '516200', '华安中证全指证券公司ETF', '0.25', '93.02', '6.12', '0.0153', 3),
    @(39, '004837', '中融鑫价值灵活配置混合C', '0.29', '63.92', '4.77', '0.0138', 6),
    @(40, '003238', '新华外延增长主题灵活配置混合', '0.42', '87.12', '3.12', '0.0131', 2),
    @(41, '012977', '瑞达鑫红量化6个月持有混合型证券投资基金A', '1.04', '94.56', '1.00', '0.0104', 7),
    @(42, '012978', '瑞达鑫红量化6个月持有混合型证券投资基金C', '0.17', '94.56', '1.00', '0.0017', 7)
)


foreach ($row in $fundData) {
    $r = $row[0] + 2

    # Column A: row index number, styled like the other sheets (bold/border)
    $newSheet.Range("A$r").Value = $row[0]
    $q4Sheet.Range("A2").Copy()
    $newSheet.Range("A$r").PasteSpecial(-4122)

    # Column B: fund code - force text so leading zeros are preserved
    $cellB = $newSheet.Range("B$r")
    $cellB.NumberFormat = "@"
    $cellB.Value = $row[1]
    $cellB.ClearFormats()

    # Column C: fund name (plain text, never numeric-looking)
    $newSheet.Range("C$r").Value = $row[2]

    # Column D: fund scale - stored as text in the source data
    $cellD = $newSheet.Range("D$r")
    $cellD.NumberFormat = "@"
    $cellD.Value = $row[3]
    $cellD.ClearFormats()

    # Column E: total stock position - stored as text in the source data
    $cellE = $newSheet.Range("E$r")
    $cellE.NumberFormat = "@"
    $cellE.Value = $row[4]
    $cellE.ClearFormats()

    # Column F: position ratio - stored as text in the source data
    $cellF = $newSheet.Range("F$r")
    $cellF.NumberFormat = "@"
    $cellF.Value = $row[5]
    $cellF.ClearFormats()

    # Column G: market value held - stored as text in the source data
    $cellG = $newSheet.Range("G$r")
    $cellG.NumberFormat = "@"
    $cellG.Value = $row[6]
    $cellG.ClearFormats()

    # Column H: position rank - actual number
    $newSheet.Range("H$r").Value = $row[7]
}

# ----------------------------------------------------------------------
# 3. Insert a new leading row into "总计" summarizing 2022-Q1, pushing the
#    existing rows down by one (their index column is recomputed below).
# ----------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 43
$totalSheet.Range("D2").Value = 62.94

Write-Host "2022-Q1 sheet added and 总计 sheet updated."
